$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..41 : columns A (debut_mvt), B (fin_mvt), C (acceleration_max)
$data = @(
    @(2, 44016.01101173611, 44016.01101173611, 18.80551484962127),
    @(3, 44016.13601262731, 44016.13601262731, 15.57982196646137),
    @(4, 44016.26101241898, 44016.26101241898, 19.46270824064371),
    @(5, 44016.38601284722, 44016.38601284722, 18.2398743872448),
    @(6, 44016.46934637731, 44016.46934637731, 16.7297456755481),
    @(7, 44016.55268054398, 44016.55268054398, 19.07880620574052),
    @(8, 44016.63601412037, 44016.63601412037, 18.9453468946629),
    @(9, 44016.71934768518, 44016.71934768518, 15.4221518268304),
    @(10, 44016.8026812963, 44016.8026812963, 14.12743757087701),
    @(11, 44016.88601451389, 44016.88601451389, 11.64253309010407),
    @(12, 44016.96934958333, 44016.96934958333, 19.10189268745386),
    @(13, 44017.05268149306, 44017.05268149306, 20.18016234122768),
    @(14, 44017.09434935185, 44017.09434935185, 17.32925357398156),
    @(15, 44017.13601577546, 44017.13601577546, 17.29736557213202),
    @(16, 44017.17768247685, 44017.17768247685, 14.78937019265054),
    @(17, 44017.21934922454, 44017.21934922454, 12.96392668670762),
    @(18, 44017.26101645833, 44017.26101645833, 17.06014353734072),
    @(19, 44017.34435050926, 44017.34435050926, 19.32721592894027),
    @(20, 44017.36518381944, 44017.36518381944, 18.28503239591587),
    @(21, 44017.38601743056, 44017.38601743056, 14.78444079222214),
    @(22, 44017.4068509838, 44017.4068509838, 13.61043070325539),
    @(23, 44017.42768402777, 44017.42768402777, 12.47544836088552),
    @(24, 44017.44851708334, 44017.44851708334, 11.88504043425403),
    @(25, 44017.46935116898, 44017.46935116898, 13.07926754730207),
    @(26, 44017.49018402777, 44017.49018402777, 18.27943959873594),
    @(27, 44017.51101752315, 44017.51101752315, 18.8105504946702),
    @(28, 44017.53185171296, 44017.53185171296, 20.93276519338268),
    @(29, 44017.55268495371, 44017.55268495371, 18.26734424055163),
    @(30, 44017.57351759259, 44017.57351759259, 14.269456528228),
    @(31, 44017.59435203703, 44017.59435203703, 15.09481957186472),
    @(32, 44017.61518556713, 44017.61518556713, 13.74438795180833),
    @(33, 44017.63601891204, 44017.63601891204, 15.31706305991288),
    @(34, 44017.65685251157, 44017.65685251157, 14.50166824687897),
    @(35, 44017.67768590278, 44017.67768590278, 18.87703936306888),
    @(36, 44017.69851924769, 44017.69851924769, 19.39444178700686),
    @(37, 44017.71935284722, 44017.71935284722, 18.22805576955546),
    @(38, 44017.74018612268, 44017.74018612268, 17.30990101423312),
    @(39, 44017.76101953704, 44017.76101953704, 13.9500873976135),
    @(40, 44017.78185297453, 44017.78185297453, 13.10761160536024),
    @(41, 44017.80268631945, 44017.80268631945, 18.6711505101691)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Row 42 no longer exists in the updated sheet - delete it entirely
$ws.Rows.Item(42).Delete()
